$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1861861861861862
    "C2" = 0.5555555555555556
    "J2" = 0.006006006006006006
    "O2" = 0.003003003003003003
    "P2" = 0.1561561561561562
    "S2" = 0.09309309309309309
    "C3" = 0.02072538860103627
    "J3" = 0.04145077720207254
    "P3" = 0.7150259067357513
    "S3" = 0.2227979274611399
    "P4" = 0.7105263157894737
    "S4" = 0.2894736842105263
    "B6" = 0.07894736842105263
    "D6" = 0.008771929824561403
    "F6" = 0.04385964912280702
    "J6" = 0.2280701754385965
    "O6" = 0.01754385964912281
    "Q6" = 0.1535087719298246
    "R6" = 0.07456140350877193
    "S6" = 0.3947368421052632
    "B7" = 0.1189189189189189
    "D7" = 0.03243243243243243
    "F7" = 0.02702702702702703
    "J7" = 0.1297297297297297
    "O7" = 0.01081081081081081
    "Q7" = 0.1351351351351351
    "R7" = 0.0918918918918919
    "S7" = 0.4540540540540541
    "B8" = 0.07822410147991543
    "D8" = 0.008456659619450317
    "F8" = 0.05496828752642706
    "J8" = 0.09090909090909091
    "O8" = 0.01902748414376321
    "Q8" = 0.1818181818181818
    "R8" = 0.08456659619450317
    "S8" = 0.4820295983086681
    "B9" = 0.1069182389937107
    "D9" = 0.01886792452830189
    "E9" = 0.006289308176100629
    "F9" = 0.06289308176100629
    "J9" = 0.119496855345912
    "O9" = 0.01257861635220126
    "Q9" = 0.1761006289308176
    "R9" = 0.1069182389937107
    "S9" = 0.389937106918239
    "B10" = 0.1250919793966151
    "D10" = 0.01766004415011038
    "F10" = 0.08462104488594555
    "J10" = 0.08462104488594555
    "O10" = 0.01103752759381899
    "Q10" = 0.2060338484179544
    "R10" = 0.08094186902133922
    "S10" = 0.3899926416482708
    "G11" = 0.1607717041800643
    "J11" = 0.1061093247588424
    "K11" = 0.2218649517684887
    "L11" = 0.4919614147909968
    "S11" = 0.01929260450160772
    "G12" = 0.7116564417177914
    "J12" = 0.1595092024539877
    "K12" = 0.0245398773006135
    "L12" = 0.05521472392638037
    "S12" = 0.049079754601227
    "G13" = 0.7027027027027027
    "J13" = 0.2432432432432433
    "S13" = 0.05405405405405406
    "F15" = 0.009302325581395349
    "H15" = 0.1813953488372093
    "I15" = 0.06511627906976744
    "J15" = 0.427906976744186
    "K15" = 0.04186046511627907
    "M15" = 0.009302325581395349
    "O15" = 0.04651162790697674
    "S15" = 0.2186046511627907
    "F16" = 0.009345794392523364
    "H16" = 0.1775700934579439
    "I16" = 0.05607476635514019
    "J16" = 0.4065420560747663
    "K16" = 0.1074766355140187
    "M16" = 0.01401869158878505
    "O16" = 0.04205607476635514
    "S16" = 0.1869158878504673
    "F17" = 0.008849557522123894
    "H17" = 0.163716814159292
    "I17" = 0.08628318584070796
    "J17" = 0.4557522123893805
    "K17" = 0.0752212389380531
    "M17" = 0.00663716814159292
    "O17" = 0.06858407079646017
    "S17" = 0.1349557522123894
    "F18" = 0.015
    "H18" = 0.16
    "I18" = 0.04
    "J18" = 0.575
    "K18" = 0.05
    "M18" = 0.025
    "O18" = 0.035
    "S18" = 0.1
    "F19" = 0.01322199025748086
    "H19" = 0.2066805845511482
    "I19" = 0.05984690327070286
    "J19" = 0.3820459290187891
    "K19" = 0.1092553931802366
    "M19" = 0.0173973556019485
    "N19" = 0.002783576896311761
    "O19" = 0.07306889352818371
    "S19" = 0.1356993736951983
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
